$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C10").Value = 8393
$ws.Range("C11:C16").Value = 7786
$ws.Range("C17:C60").Value = 7685
$ws.Range("C61:C85").Value = 7660
$ws.Range("C86:C91").Value = 7318
